# "Fruta / hortaliza, semanal" - weekly update:
# a new weekly price record is inserted as row 19 (Feria Lagunitas de
# Puerto Montt - Pomelo - Start Ruby - Primera), pushing all the
# subsequent rows (old rows 19-141) down by one (new rows 20-142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 19; this shifts rows
# 19..141 down to 20..142 and keeps the column D date formatting.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with this week's data.
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44462
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100102
$ws.Range("H19").Value = "Cítricos"
$ws.Range("I19").Value = 100102006
$ws.Range("J19").Value = "Pomelo"
$ws.Range("K19").Value = "Start Ruby"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 12000
$ws.Range("O19").Value = 12000
$ws.Range("P19").Value = 12000
$ws.Range("Q19").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 857
$ws.Range("T19").Value = 14
